$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells, copying the style of the existing header row (A1:F1)
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Update existing values in row 2
$ws.Range("B2").Value = 0.08557473815421481
$ws.Range("C2").Value = 0.9991882852584594
$ws.Range("D2").Value = 0.2201936378481086

# Add new values in row 2
$ws.Range("G2").Value = 0.1180509527000443
$ws.Range("H2").Value = 0.991
